$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue 2 4 '56.624.20'
Set-TextValue 2 5 '  +4.40%  '
Set-TextValue 3 4 '3.010.09'
Set-TextValue 3 5 '  +5.09%  '
Set-TextValue 4 5 '  +0.05%  '
Set-TextValue 5 4 '509.92'
Set-TextValue 5 5 '  +8.35%  '
Set-TextValue 6 4 '137.07'
Set-TextValue 6 5 '  +9.02%  '
Set-TextValue 7 5 '  +0.01%  '
Set-TextValue 8 4 '0.433'
Set-TextValue 8 5 '  +8.06%  '
Set-TextValue 9 4 '7.60'
Set-TextValue 9 5 '  +14.65%  '
Set-TextValue 10 5 '  +13.57%  '
Set-TextValue 11 4 '0.354'
Set-TextValue 11 5 '  +7.07%  '
Set-TextValue 12 5 '  +4.50%  '
Set-TextValue 13 4 '3.520.61'
Set-TextValue 13 5 '  +5.05%  '
Set-TextValue 14 4 '25.76'
Set-TextValue 14 5 '  +10.53%  '
Set-TextValue 15 5 '  +15.14%  '
Set-TextValue 16 4 '56.657.32'
Set-TextValue 16 5 '  +4.60%  '
Set-TextValue 17 4 '3.005.74'
Set-TextValue 17 5 '  +5.11%  '
Set-TextValue 18 4 '5.84'
Set-TextValue 18 5 '  +9.26%  '
Set-TextValue 19 4 '12.50'
Set-TextValue 19 5 '  +9.37%  '
Set-TextValue 20 4 '7.88'
Set-TextValue 20 5 '  +11.82%  '
Set-TextValue 21 4 '326.78'
Set-TextValue 21 5 '  +11.30%  '
Set-TextValue 22 4 '1.00'
Set-TextValue 22 5 '  -0.11%  '
Set-TextValue 23 4 '0.477'
Set-TextValue 24 4 '62.54'
Set-TextValue 24 5 '  +6.42%  '
Set-TextValue 25 4 '1.00'
Set-TextValue 25 5 '  -0.29%  '
Set-TextValue 26 4 '0.165'
Set-TextValue 26 5 '  +7.30%  '
Set-TextValue 27 4 '0.0₃0918'
Set-TextValue 27 5 '  +14.06%  '
Set-TextValue 28 4 '6.54'
Set-TextValue 28 5 '  +5.85%  '
Set-TextValue 29 4 '6.96'
Set-TextValue 29 5 '  +12.09%  '
Set-TextValue 30 4 '1.24'
Set-TextValue 30 5 '  +9.78%  '
Set-TextValue 31 4 '1.77'
Set-TextValue 31 5 '  +9.89%  '
Set-TextValue 32 4 '20.62'
Set-TextValue 32 5 '  +9.66%  '
Set-TextValue 33 4 '156.57'
Set-TextValue 33 5 '  +15.69%  '
Set-TextValue 34 4 '4.51'
Set-TextValue 34 5 '  +6.91%  '
Set-TextValue 35 4 '5.62'
Set-TextValue 35 5 '  +3.66%  '
Set-TextValue 36 4 '1.27'
Set-TextValue 36 5 '  +3.92%  '
Set-TextValue 37 4 '0.0675'
Set-TextValue 37 5 '  +9.57%  '
Set-TextValue 38 4 '23.80'
Set-TextValue 38 5 '  +2.92%  '
Set-TextValue 39 4 '3.042.51'
Set-TextValue 39 5 '  +5.44%  '
Set-TextValue 40 4 '36.53'
Set-TextValue 40 5 '  +4.37%  '
Set-TextValue 41 5 '  +0.08%  '
Set-TextValue 42 4 '0.648'
Set-TextValue 42 5 '  +7.75%  '
Set-TextValue 43 4 '2.266.12'
Set-TextValue 43 5 '  +11.07%  '
Set-TextValue 44 2 'ONDO'
Set-TextValue 44 3 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 44 4 '0.999'
Set-TextValue 44 5 '  +5.33%  '
Set-TextValue 45 2 'Stacks'
Set-TextValue 45 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 45 4 '1.41'
Set-TextValue 45 5 '  +7.38%  '
Set-TextValue 46 4 '3.61'
Set-TextValue 46 5 '  +6.38%  '
Set-TextValue 47 4 '1.98'
Set-TextValue 47 5 '  +23.97%  '
Set-TextValue 48 4 '0.0237'
Set-TextValue 48 5 '  +11.86%  '
Set-TextValue 49 4 '5.79'
Set-TextValue 49 5 '  +8.66%  '
Set-TextValue 50 4 '19.14'
Set-TextValue 50 5 '  +7.01%  '
Set-TextValue 51 4 '0.0874'
Set-TextValue 51 5 '  +10.92%  '
